$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing Basket rows: correct the "api/vi/" typo to "api/v1/"
# and swap Update/Checkout use-case text to match corrected rows (row 8 & 9)
$ws.Range("C8").Value = "api/v1/Basket"
$ws.Range("D8").Value = "Get Basket and Items with Username"

$ws.Range("C9").Value = "api/v1/Basket"
$ws.Range("D9").Value = "Update Basket and Items (add - remove item on basket)"

$ws.Range("C10").Value = "api/v1/Basket/{id}"
$ws.Range("D10").Value = "Delete Basket"

$ws.Range("C11").Value = "api/v1/Basket/Checkout"
$ws.Range("D11").Value = "Checkout Basket"

# Add new "Order" rows (12 and 13)
$ws.Range("A12").Value = "Order"
$ws.Range("B12").Value = "GET"
$ws.Range("C12").Value = "api/v1/Order"
$ws.Range("D12").Value = "Get Orders with username"

$ws.Range("A13").Value = "Order"
$ws.Range("B13").Value = "PUT"
$ws.Range("C13").Value = "api/v1/Order/"
$ws.Range("D13").Value = "Post Orders with username (for testing)"

$ws.Range("C13").Select()
